$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: new equipment-log entry (this row previously held just the "**" marker in A7)
$ws.Range("A7").Value = "31/01/2018"

# "3012" is pure-digit text; a plain .Value assignment would coerce it to a number,
# so force it to text via a formula that evaluates to a string, then freeze the result
# as a value (Copy + PasteSpecial values-only) so it lands as a shared string with no
# extra number-format/quote-prefix style.
$ws.Range("B7").Formula = "=""3012"""
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)

$ws.Range("C7").Value = "Зазубрини в місті відрізу контакту"

$ws.Range("D7").Formula = "=""6"""
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)

# Row 8: the "**" marker moves down here
$ws.Range("A8").Value = "**"

$excel.CutCopyMode = $false
